{"js": "// Split the single \"Bibliografia\" paragraph's run into multiple lines by\n// inserting a manual line break (serialized as w:br) right before each\n// \"[n]\" citation marker (n = 2..5). We go from the last marker to the\n// first so inserting text for an earlier marker never moves the text we\n// still need to search for.\nconst body = context.document.body;\n\nconst markers = [\n  \"[5] Diretrizes\",\n  \"[4] Kaul\",\n  \"[3] Mueller\",\n  \"[2] Zachary\",\n];\n\nfor (const marker of markers) {\n  const results = body.search(marker, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    // \"\\v\" (vertical tab) is Word's manual line-break character; Office.js\n    // serializes it as a <w:br/> between two <w:t> runs instead of\n    // starting a brand-new paragraph.\n    results.items[0].insertText(\"\\v\", Word.InsertLocation.before);\n    await context.sync();\n  }\n}\n", "ps1": "# Split the single \"Bibliografia\" paragraph's run into multiple lines by\n# inserting a manual line break (w:br) right before each \"[n]\" citation\n# marker (n = 2..5), going from the last marker to the first so earlier\n# inserts don't shift the character offsets of markers still to be found.\n\n$d = $word.ActiveDocument\n\n$markers = @(\n    \"[5] Diretrizes\",\n    \"[4] Kaul\",\n    \"[3] Mueller\",\n    \"[2] Zachary\"\n)\n\nforeach ($marker in $markers) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $found = $rng.Find.Execute($marker)\n    if ($found) {\n        $rng.Collapse(1)              # wdCollapseStart\n        $rng.InsertBefore([char]11)   # vertical-tab -> manual line break (w:br)\n    }\n}\n\n$d.Save()\n"}
